$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the source URL column
$ws.Range("S1").Value = "Quelle_URL"

# Copy the header formatting (bold, centered, bordered) from the existing
# "Quelle" header cell so the new column matches the rest of the header row
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

# Column width for the new column S (stored OOXML width of 45 chars)
$ws.Columns.Item(19).ColumnWidth = 44.14

# Source URLs per row
$ws.Range("S2").Value = "https://www.refurbed.de/dell-latitude-3550"
$ws.Range("S3").Value = "https://www.refurbed.de/dell-precision-5550"
$ws.Range("S4").Value = "https://www.kleinanzeigen.de"
$ws.Range("S5").Value = "https://www.refurbed.de"
$ws.Range("S6").Value = "https://www.refurbed.de"
$ws.Range("S7").Value = "https://www.refurbed.de"
$ws.Range("S8").Value = "https://www.refurbed.de"
$ws.Range("S9").Value = "https://www.refurbed.de"
$ws.Range("S10").Value = "https://www.refurbed.de"
$ws.Range("S11").Value = "https://www.refurbed.de"
$ws.Range("S12").Value = "https://www.orbit365.de"
